$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.053279369324243
$ws.Range("D2").Value = 1.052026283836328
$ws.Range("E2").Value = 1.056960664408924
$ws.Range("F2").Value = 1.05224569416336
$ws.Range("I2").Value = 1.047820075865529
$ws.Range("J2").Value = 1.058297547625609
$ws.Range("K2").Value = 1.054776156228006
$ws.Range("L2").Value = 1.059696942303818
$ws.Range("M2").Value = 1.054994959132021
$ws.Range("N2").Value = 1.059800451531663

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.054947396615801
$ws.Range("D3").Value = 1.053336482996229
$ws.Range("E3").Value = 1.058579903389779
$ws.Range("F3").Value = 1.054539013005928
$ws.Range("I3").Value = 1.048467183520795
$ws.Range("J3").Value = 1.059612907994027
$ws.Range("K3").Value = 1.055897841636719
$ws.Range("L3").Value = 1.061127876392949
$ws.Range("M3").Value = 1.057097289088976
$ws.Range("N3").Value = 1.061117679862678

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.056021723291454
$ws.Range("D4").Value = 1.05417937446767
$ws.Range("E4").Value = 1.059623026920053
$ws.Range("F4").Value = 1.056017440773013
$ws.Range("I4").Value = 1.048881772013382
$ws.Range("J4").Value = 1.060458718591526
$ws.Range("K4").Value = 1.05661821832878
$ws.Range("L4").Value = 1.062048683545006
$ws.Range("M4").Value = 1.058451815580914
$ws.Range("N4").Value = 1.061964691608241

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.05647219489016
$ws.Range("D5").Value = 1.054532571193013
$ws.Range("E5").Value = 1.060060467144906
$ws.Range("F5").Value = 1.056637687110602
$ws.Range("I5").Value = 1.049055085567398
$ws.Range("J5").Value = 1.060813042949542
$ws.Range("K5").Value = 1.056919781235104
$ws.Range("L5").Value = 1.062434588151751
$ws.Range("M5").Value = 1.059019895295478
$ws.Range("N5").Value = 1.062319519147487

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.056547762745771
$ws.Range("D6").Value = 1.054591807246763
$ws.Range("E6").Value = 1.060133852008261
$ws.Range("F6").Value = 1.056741754807033
$ws.Range("I6").Value = 1.049084128531342
$ws.Range("J6").Value = 1.060872462530907
$ws.Range("K6").Value = 1.056970340178333
$ws.Range("L6").Value = 1.062499313250278
$ws.Range("M6").Value = 1.059115199317159
$ws.Range("N6").Value = 1.062379023111472

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.056027747104445
$ws.Range("D7").Value = 1.054184098414323
$ws.Range("E7").Value = 1.059628876269231
$ws.Range("F7").Value = 1.056025733538908
$ws.Range("I7").Value = 1.048884091672456
$ws.Range("J7").Value = 1.060463457995429
$ws.Range("K7").Value = 1.056622252848274
$ws.Range("L7").Value = 1.062053844721526
$ws.Range("M7").Value = 1.05845941160538
$ws.Range("N7").Value = 1.061969437742641

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.053844134847235
$ws.Range("D8").Value = 1.052470095630483
$ws.Range("E8").Value = 1.05750886513508
$ws.Range("F8").Value = 1.053021891216652
$ws.Range("I8").Value = 1.048039630903881
$ws.Range("J8").Value = 1.058743191415595
$ws.Range("K8").Value = 1.055156368943634
$ws.Range("L8").Value = 1.060181601227694
$ws.Range("M8").Value = 1.05570667523497
$ws.Range("N8").Value = 1.060246728186947

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.049957082966527
$ws.Range("D9").Value = 1.049411542713679
$ws.Range("E9").Value = 1.053736746421368
$ws.Range("F9").Value = 1.047685111514485
$ws.Range("I9").Value = 1.046519467990766
$ws.Range("J9").Value = 1.055670330180893
$ws.Range("K9").Value = 1.05253098216545
$ws.Range("L9").Value = 1.056842529889278
$ws.Range("M9").Value = 1.050810036274727
$ws.Range("N9").Value = 1.057169503136829

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.047337965441378
$ws.Range("D10").Value = 1.047345694573819
$ws.Range("E10").Value = 1.05119625190482
$ws.Range("F10").Value = 1.044095841220334
$ws.Range("I10").Value = 1.045483786890334
$ws.Range("J10").Value = 1.053592673785224
$ws.Range("K10").Value = 1.050751240131627
$ws.Range("L10").Value = 1.054588420433054
$ws.Range("M10").Value = 1.047512761511578
$ws.Range("N10").Value = 1.055088896230768

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.04619697054639
$ws.Range("D11").Value = 1.046444554816267
$ws.Range("E11").Value = 1.050089793714729
$ws.Range("F11").Value = 1.042533727177106
$ws.Range("I11").Value = 1.045029904949248
$ws.Range("J11").Value = 1.052685866893463
$ws.Range("K11").Value = 1.049973365679109
$ws.Range("L11").Value = 1.053605438108778
$ws.Range("M11").Value = 1.046076776925966
$ws.Range("N11").Value = 1.054180801569208

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.04577209101857
$ws.Range("D12").Value = 1.046108816400454
$ws.Range("E12").Value = 1.049677818151629
$ws.Range("F12").Value = 1.041952255478086
$ws.Range("I12").Value = 1.044860485205822
$ws.Range("J12").Value = 1.052347938842237
$ws.Range("K12").Value = 1.049683321279959
$ws.Range("L12").Value = 1.053239248673391
$ws.Range("M12").Value = 1.045542112006216
$ws.Range("N12").Value = 1.053842393621389

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.04586327766443
$ws.Range("D13").Value = 1.046180879706703
$ws.Range("E13").Value = 1.049766233412417
$ws.Range("F13").Value = 1.042077039586424
$ws.Range("I13").Value = 1.044896863981004
$ws.Range("J13").Value = 1.052420475683613
$ws.Range("K13").Value = 1.049745587220458
$ws.Range("L13").Value = 1.053317846174657
$ws.Range("M13").Value = 1.045656857817304
$ws.Range("N13").Value = 1.053915033473399

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.046161871716986
$ws.Range("D14").Value = 1.046416823404127
$ws.Range("E14").Value = 1.050055759982164
$ws.Range("F14").Value = 1.042485687935262
$ws.Range("I14").Value = 1.045015917621174
$ws.Range("J14").Value = 1.052657956219363
$ws.Range("K14").Value = 1.049949413223649
$ws.Range("L14").Value = 1.053575190660877
$ws.Range("M14").Value = 1.046032607614767
$ws.Range("N14").Value = 1.05415285125875

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.046345703676871
$ws.Range("D15").Value = 1.046562060963666
$ws.Range("E15").Value = 1.05023401521609
$ws.Range("F15").Value = 1.04273730505936
$ws.Range("I15").Value = 1.045089160425265
$ws.Range("J15").Value = 1.05280412942292
$ws.Range("K15").Value = 1.050074849824401
$ws.Range("L15").Value = 1.053733607117855
$ws.Range("M15").Value = 1.046263949198247
$ws.Range("N15").Value = 1.054299232045022

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.047413541883532
$ws.Range("D16").Value = 1.047405359092497
$ws.Range("E16").Value = 1.051269546763103
$ws.Range("F16").Value = 1.044199342921196
$ws.Range("I16").Value = 1.045513794157531
$ws.Range("J16").Value = 1.053652702623073
$ws.Range("K16").Value = 1.050802710950954
$ws.Range("L16").Value = 1.054653509544241
$ws.Range("M16").Value = 1.047607886189938
$ws.Range("N16").Value = 1.055149010316452

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.048081502671864
$ws.Range("D17").Value = 1.047932552059543
$ws.Range("E17").Value = 1.051917375334676
$ws.Range("F17").Value = 1.04511428738047
$ws.Range("I17").Value = 1.045778694296515
$ws.Range("J17").Value = 1.054183055203046
$ws.Range("K17").Value = 1.051257327919152
$ws.Range("L17").Value = 1.055228665930372
$ws.Range("M17").Value = 1.048448668943242
$ws.Range("N17").Value = 1.055680116057921

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.048470448994689
$ws.Range("D18").Value = 1.048239418112302
$ws.Range("E18").Value = 1.052294626284493
$ws.Range("F18").Value = 1.045647196042311
$ws.Range("I18").Value = 1.045932683482693
$ws.Range("J18").Value = 1.054491710836963
$ws.Range("K18").Value = 1.051521801677933
$ws.Range("L18").Value = 1.055563477115891
$ws.Range("M18").Value = 1.048938290573691
$ws.Range("N18").Value = 1.055989210018242

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.04860295782891
$ws.Range("D19").Value = 1.048343944153381
$ws.Range("E19").Value = 1.052423155294548
$ws.Range("F19").Value = 1.045828775933254
$ws.Range("I19").Value = 1.045985101557959
$ws.Range("J19").Value = 1.054596838067895
$ws.Range("K19").Value = 1.051611862844494
$ws.Range("L19").Value = 1.055677526517265
$ws.Range("M19").Value = 1.049105105575064
$ws.Range("N19").Value = 1.0560944865419

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.048009905662987
$ws.Range("D20").Value = 1.047876055252284
$ws.Range("E20").Value = 1.051847933407566
$ws.Range("F20").Value = 1.045016201707927
$ws.Range("I20").Value = 1.045750327163232
$ws.Range("J20").Value = 1.054126224929087
$ws.Range("K20").Value = 1.051208624022414
$ws.Range("L20").Value = 1.055167026337705
$ws.Range("M20").Value = 1.048358543178751
$ws.Range("N20").Value = 1.055623205078455

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.046073972776752
$ws.Range("D21").Value = 1.046347372057573
$ws.Range("E21").Value = 1.049970529093324
$ws.Range("F21").Value = 1.042365385548131
$ws.Range("I21").Value = 1.044980882245399
$ws.Range("J21").Value = 1.052588054658932
$ws.Range("K21").Value = 1.049889422274584
$ws.Range("L21").Value = 1.053499438743547
$ws.Range("M21").Value = 1.045921994204613
$ws.Range("N21").Value = 1.054082850430085

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.044850609968144
$ws.Range("D22").Value = 1.045380345042396
$ws.Range("E22").Value = 1.048784403076564
$ws.Range("F22").Value = 1.040691557894617
$ws.Range("I22").Value = 1.044492305949719
$ws.Range("J22").Value = 1.051614572526308
$ws.Range("K22").Value = 1.049053571602112
$ws.Range("L22").Value = 1.05244478049386
$ws.Range("M22").Value = 1.044382635602963
$ws.Range("N22").Value = 1.053107985841173

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.045499731105658
$ws.Range("D23").Value = 1.045893549232011
$ws.Range("E23").Value = 1.049413742299693
$ws.Range("F23").Value = 1.041579577496662
$ws.Range("I23").Value = 1.044751768366531
$ws.Range("J23").Value = 1.052131245624276
$ws.Range("K23").Value = 1.049497286903332
$ws.Range("L23").Value = 1.05300446855802
$ws.Range("M23").Value = 1.045199393315448
$ws.Range("N23").Value = 1.053625392674203

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.048042259320016
$ws.Range("D24").Value = 1.047901585694271
$ws.Range("E24").Value = 1.051879313129491
$ws.Range("F24").Value = 1.045060524757272
$ws.Range("I24").Value = 1.04576314666405
$ws.Range("J24").Value = 1.054151906215143
$ws.Range("K24").Value = 1.051230633368788
$ws.Range("L24").Value = 1.055194880679103
$ws.Range("M24").Value = 1.048399269580356
$ws.Range("N24").Value = 1.055648922834882

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.050966769016584
$ws.Range("D25").Value = 1.050206897071066
$ws.Range("E25").Value = 1.05471637261193
$ws.Range("F25").Value = 1.049070179606773
$ws.Range("I25").Value = 1.046916338292333
$ws.Range("J25").Value = 1.056469778352333
$ws.Range("K25").Value = 1.053214827420182
$ws.Range("L25").Value = 1.05771061616944
$ws.Range("M25").Value = 1.052081580282259
$ws.Range("N25").Value = 1.057970086616371
